$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Type" header to "Type * " (now a required field, like the other
# asterisked headers in the sheet).
$ws.Range("C1").Value = "Type * "

# Add a new "Commitment Date" column so FX conversions can be pinned to a date.
$ws.Range("J1").Value = "Commitment Date"

# Commitment date for each existing row: 1/20/2023 (serial 44946).
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Value = 44946
    $cell.NumberFormat = "mm-dd-yy"
}

# Widen the new column to fit the date values, matching the other data columns.
$ws.Columns.Item(10).ColumnWidth = 15.625

# Move the selection to the newly added column, as Excel would after data entry.
$ws.Range("J3:J7").Select
